$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("CustomerID") before the existing Status column,
# shifting Status/Description one column to the right (C, D).
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "CustomerID"

$custIds = @("ANATR","TORTU","VAFFE","MEREP","PRINI","HILAA","LEHMS","QUEEN","LAMAI","FRANK","WARTH","BONAP","PERIC","MEREP","QUICK")

for ($i = 0; $i -lt $custIds.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $custIds[$i]
}

$dataRange = $ws.Range("B2:B16")
$dataRange.NumberFormat = "#"
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 11
$dataRange.Font.Color = 0

[void]$ws.Range("B16").Select()
